# Updates cryptos list data (prices, 1h volume %, and one coin rename)
# matching the commit 'Updated cryptos list on Wed Jan 31 21:09:59 UTC 2024 with GitHub Actions'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$value) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (prevents numeric/percentage auto-conversion of strings that look
    # like numbers, e.g. "2.83" or "  -2.20%  "), then resetting the style
    # back to Normal drops the quote-prefix flag Excel would otherwise
    # stamp on the cell, keeping formatting identical to the original file.
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell 'D2' '42.593.10'
Set-TextCell 'E2' '  -2.20%  '
Set-TextCell 'D3' '2.286.33'
Set-TextCell 'E3' '  -3.65%  '
Set-TextCell 'E4' '  -0.04%  '
Set-TextCell 'D5' '300.61'
Set-TextCell 'E5' '  -2.89%  '
Set-TextCell 'D6' '98.08'
Set-TextCell 'E6' '  -6.51%  '
Set-TextCell 'E7' '  -0.58%  '
Set-TextCell 'E8' '  -0.02%  '
Set-TextCell 'D9' '0.502'
Set-TextCell 'E9' '  -3.55%  '
Set-TextCell 'D10' '34.50'
Set-TextCell 'E10' '  -4.21%  '
Set-TextCell 'D11' '51.01'
Set-TextCell 'E11' '  -4.57%  '
Set-TextCell 'E12' '  -2.54%  '
Set-TextCell 'E13' '  +0.28%  '
Set-TextCell 'D14' '6.72'
Set-TextCell 'E14' '  -4.02%  '
Set-TextCell 'D15' '2.640.76'
Set-TextCell 'E15' '  -3.77%  '
Set-TextCell 'D16' '15.45'
Set-TextCell 'E16' '  -1.39%  '
Set-TextCell 'D17' '2.291.22'
Set-TextCell 'E17' '  -3.50%  '
Set-TextCell 'D18' '0.792'
Set-TextCell 'E18' '  -2.29%  '
Set-TextCell 'D19' '42.506.29'
Set-TextCell 'E19' '  -2.33%  '
Set-TextCell 'D20' '11.60'
Set-TextCell 'E20' '  -2.73%  '
Set-TextCell 'D21' '0.0₃0896'
Set-TextCell 'E21' '  -2.32%  '
Set-TextCell 'D22' '6.02'
Set-TextCell 'E22' '  -4.51%  '
Set-TextCell 'D23' '66.97'
Set-TextCell 'E23' '  -2.03%  '
Set-TextCell 'D24' '234.99'
Set-TextCell 'E24' '  -2.47%  '
Set-TextCell 'E25' '  -4.99%  '
Set-TextCell 'D26' '2.50'
Set-TextCell 'E26' '  -4.39%  '
Set-TextCell 'E27' '  -0.30%  '
Set-TextCell 'D28' '24.59'
Set-TextCell 'E28' '  -4.64%  '
Set-TextCell 'E29' '  -1.01%  '
Set-TextCell 'D30' '34.19'
Set-TextCell 'E30' '  -6.52%  '
Set-TextCell 'D31' '164.25'
Set-TextCell 'E31' '  +2.11%  '
Set-TextCell 'D32' '9.12'
Set-TextCell 'E32' '  -4.27%  '
Set-TextCell 'E33' '  +0.07%  '
Set-TextCell 'E34' '  -5.08%  '
Set-TextCell 'E35' '  -4.98%  '
Set-TextCell 'D36' '0.0699'
Set-TextCell 'E36' '  -5.60%  '
Set-TextCell 'D37' '4.35'
Set-TextCell 'E37' '  -6.55%  '
Set-TextCell 'D38' '2.83'
Set-TextCell 'E38' '  -8.69%  '
Set-TextCell 'D39' '16.23'
Set-TextCell 'E39' '  -11.26%  '
Set-TextCell 'D40' '1.78'
Set-TextCell 'E40' '  -7.97%  '
Set-TextCell 'E41' '  -4.58%  '
Set-TextCell 'E42' '  -2.94%  '
Set-TextCell 'E43' '  -9.33%  '
Set-TextCell 'D44' '1.966.77'
Set-TextCell 'E44' '  -3.17%  '
Set-TextCell 'D45' '0.0283'
Set-TextCell 'E45' '  -2.82%  '
Set-TextCell 'D46' '18.16'
Set-TextCell 'E46' '  -7.98%  '
Set-TextCell 'D47' '9.76'
Set-TextCell 'E47' '  -7.56%  '
Set-TextCell 'D48' '2.86'
Set-TextCell 'E48' '  -9.10%  '
Set-TextCell 'D49' '4.72'
Set-TextCell 'E49' '  -0.37%  '
Set-TextCell 'B50' 'HuobiToken'
Set-TextCell 'C50' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D50' '2.83'
Set-TextCell 'E50' '  -4.60%  '
Set-TextCell 'D51' '2.514.79'
Set-TextCell 'E51' '  -3.31%  '
